# Insert a new weekly price record at row 88 ("Fruta / hortaliza, semanal"):
# the sheet stores the most recent week first, so the newest observation
# (week of 2023-10-12, serial 45211) is inserted above the existing
# historical rows, which all shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 88:106 down to 89:107, leaving a blank row 88 to populate.
$ws.Rows(88).Insert()

$ws.Cells.Item(88, 1).Value  = 10
$ws.Cells.Item(88, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(88, 3).Value  = "La Araucanía"
$ws.Cells.Item(88, 4).Value  = 45211
$ws.Cells.Item(88, 5).Value  = 9
$ws.Cells.Item(88, 6).Value  = 100112042
$ws.Cells.Item(88, 7).Value  = "Locoto"
$ws.Cells.Item(88, 8).Value  = "Sin especificar"
$ws.Cells.Item(88, 9).Value  = "Primera"
$ws.Cells.Item(88, 10).Value = 80
$ws.Cells.Item(88, 11).Value = 2200
$ws.Cells.Item(88, 12).Value = 2200
$ws.Cells.Item(88, 13).Value = 2200
$ws.Cells.Item(88, 14).Value = "`$/kilo"
$ws.Cells.Item(88, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(88, 16).Value = 2200
$ws.Cells.Item(88, 17).Value = 1
$ws.Cells.Item(88, 18).Value = "Hortaliza"

# Match the date-time number format used by the rest of column D.
$ws.Cells.Item(88, 4).NumberFormat = $ws.Cells.Item(89, 4).NumberFormat
